$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.401.04"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.062.92"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'592.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").Value = "3.063.54"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "'36.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "3.568.12"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "63.336.64"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "3.063.52"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'487.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("D25").Value = "'82.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").Value = "'10.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.06%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'2.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "'27.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0825"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "'3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").Value = "'2.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").Value = "'50.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "'439.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "'0.114"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.36%  "
$ws.Range("D45").Value = "'0.0365"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "2.851.87"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").Value = "'38.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "'130.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D50").Value = "'25.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").Value = "  -0.38%  "
